$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing rows 2-13: only cells that changed ----
# Row 2: columns D,L,M,N,O,P,R,S
$ws.Cells.Item(2,4).Value = 45091
$ws.Cells.Item(2,12).Value = 'Primera'
$ws.Cells.Item(2,13).Value = 220
$ws.Cells.Item(2,14).Value = 18000
$ws.Cells.Item(2,15).Value = 19000
$ws.Cells.Item(2,16).Value = 18455
$ws.Cells.Item(2,18).Value = 'Provincia de Curicó'
$ws.Cells.Item(2,19).Value = 1025

# Row 3: columns D,L,M,N,O,P,Q,R,S,T
$ws.Cells.Item(3,4).Value = 45091
$ws.Cells.Item(3,12).Value = 'Segunda'
$ws.Cells.Item(3,13).Value = 150
$ws.Cells.Item(3,14).Value = 15000
$ws.Cells.Item(3,15).Value = 15000
$ws.Cells.Item(3,16).Value = 15000
$ws.Cells.Item(3,17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(3,18).Value = 'Provincia de Curicó'
$ws.Cells.Item(3,19).Value = 833
$ws.Cells.Item(3,20).Value = 18

# Row 4: columns D,L,M,N,O,P,Q,R,S,T
$ws.Cells.Item(4,4).Value = 44708
$ws.Cells.Item(4,12).Value = 'Primera'
$ws.Cells.Item(4,13).Value = 70
$ws.Cells.Item(4,14).Value = 12000
$ws.Cells.Item(4,15).Value = 13000
$ws.Cells.Item(4,16).Value = 12571
$ws.Cells.Item(4,17).Value = '$/caja 12 kilos empedrada'
$ws.Cells.Item(4,18).Value = 'Provincia de Curicó'
$ws.Cells.Item(4,19).Value = 1048
$ws.Cells.Item(4,20).Value = 12

# Row 5: columns D,M,O,P,Q,R,S,T
$ws.Cells.Item(5,4).Value = 45077
$ws.Cells.Item(5,13).Value = 140
$ws.Cells.Item(5,15).Value = 14000
$ws.Cells.Item(5,16).Value = 12857
$ws.Cells.Item(5,17).Value = '$/caja 12 kilos granel'
$ws.Cells.Item(5,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(5,19).Value = 12857
$ws.Cells.Item(5,20).Value = 1

# Row 6: columns D,L,M,O,P,S
$ws.Cells.Item(6,4).Value = 45077
$ws.Cells.Item(6,12).Value = 'Segunda'
$ws.Cells.Item(6,13).Value = 80
$ws.Cells.Item(6,15).Value = 11000
$ws.Cells.Item(6,16).Value = 11000
$ws.Cells.Item(6,19).Value = 11000

# Row 7: columns D,L,N,O,P,R,S
$ws.Cells.Item(7,4).Value = 44742
$ws.Cells.Item(7,12).Value = 'Segunda'
$ws.Cells.Item(7,14).Value = 14000
$ws.Cells.Item(7,15).Value = 15000
$ws.Cells.Item(7,16).Value = 14500
$ws.Cells.Item(7,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(7,19).Value = 806

# Row 8: columns D,M,N,O,P,R,S
$ws.Cells.Item(8,4).Value = 44714
$ws.Cells.Item(8,13).Value = 100
$ws.Cells.Item(8,14).Value = 14000
$ws.Cells.Item(8,15).Value = 15000
$ws.Cells.Item(8,16).Value = 14500
$ws.Cells.Item(8,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(8,19).Value = 806

# Row 9: columns D,L,M,N,O,P,Q,R,S,T
$ws.Cells.Item(9,4).Value = 44334
$ws.Cells.Item(9,12).Value = 'Primera'
$ws.Cells.Item(9,13).Value = 100
$ws.Cells.Item(9,14).Value = 11000
$ws.Cells.Item(9,15).Value = 12000
$ws.Cells.Item(9,16).Value = 11500
$ws.Cells.Item(9,17).Value = '$/caja 12 kilos granel'
$ws.Cells.Item(9,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(9,19).Value = 11500
$ws.Cells.Item(9,20).Value = 1

# Row 10: columns D,M,N,O,P,Q,R,S,T
$ws.Cells.Item(10,4).Value = 45084
$ws.Cells.Item(10,13).Value = 100
$ws.Cells.Item(10,14).Value = 17000
$ws.Cells.Item(10,15).Value = 18000
$ws.Cells.Item(10,16).Value = 17500
$ws.Cells.Item(10,17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(10,18).Value = 'Región del Maule'
$ws.Cells.Item(10,19).Value = 972
$ws.Cells.Item(10,20).Value = 18

# Row 11: columns D,M,N,O,P,R,S
$ws.Cells.Item(11,4).Value = 45090
$ws.Cells.Item(11,13).Value = 150
$ws.Cells.Item(11,14).Value = 17000
$ws.Cells.Item(11,15).Value = 18000
$ws.Cells.Item(11,16).Value = 17533
$ws.Cells.Item(11,18).Value = 'Región del Maule'
$ws.Cells.Item(11,19).Value = 974

# Row 12: columns D,L,M,P,S
$ws.Cells.Item(12,4).Value = 45090
$ws.Cells.Item(12,12).Value = 'Segunda'
$ws.Cells.Item(12,13).Value = 130
$ws.Cells.Item(12,16).Value = 14462
$ws.Cells.Item(12,19).Value = 803

# Row 13: columns D,M,N,O,P,Q,R,S,T
$ws.Cells.Item(13,4).Value = 44707
$ws.Cells.Item(13,13).Value = 60
$ws.Cells.Item(13,14).Value = 12000
$ws.Cells.Item(13,15).Value = 13000
$ws.Cells.Item(13,16).Value = 12500
$ws.Cells.Item(13,17).Value = '$/caja 12 kilos empedrada'
$ws.Cells.Item(13,18).Value = 'Provincia de Curicó'
$ws.Cells.Item(13,19).Value = 1042
$ws.Cells.Item(13,20).Value = 12

# ---- Append new rows 14-15 (full rows) ----
# Row 14
$ws.Cells.Item(14,1).Value = 11
$ws.Cells.Item(14,2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(14,3).Value = 'Bíobío'
$ws.Cells.Item(14,4).Value = 44330
$ws.Cells.Item(14,5).Value = 8
$ws.Cells.Item(14,6).Value = 'Fruta'
$ws.Cells.Item(14,7).Value = 100107
$ws.Cells.Item(14,8).Value = 'Otros'
$ws.Cells.Item(14,9).Value = 100107001
$ws.Cells.Item(14,10).Value = 'Caqui'
$ws.Cells.Item(14,11).Value = 'Mankaki'
$ws.Cells.Item(14,12).Value = 'Primera'
$ws.Cells.Item(14,13).Value = 100
$ws.Cells.Item(14,14).Value = 15000
$ws.Cells.Item(14,15).Value = 16000
$ws.Cells.Item(14,16).Value = 15500
$ws.Cells.Item(14,17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(14,18).Value = 'Provincia de Curicó'
$ws.Cells.Item(14,19).Value = 861
$ws.Cells.Item(14,20).Value = 18
$ws.Cells.Item(14,4).NumberFormat = $ws.Cells.Item(13,4).NumberFormat

# Row 15
$ws.Cells.Item(15,1).Value = 11
$ws.Cells.Item(15,2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(15,3).Value = 'Bíobío'
$ws.Cells.Item(15,4).Value = 44719
$ws.Cells.Item(15,5).Value = 8
$ws.Cells.Item(15,6).Value = 'Fruta'
$ws.Cells.Item(15,7).Value = 100107
$ws.Cells.Item(15,8).Value = 'Otros'
$ws.Cells.Item(15,9).Value = 100107001
$ws.Cells.Item(15,10).Value = 'Caqui'
$ws.Cells.Item(15,11).Value = 'Mankaki'
$ws.Cells.Item(15,12).Value = 'Primera'
$ws.Cells.Item(15,13).Value = 50
$ws.Cells.Item(15,14).Value = 14000
$ws.Cells.Item(15,15).Value = 15000
$ws.Cells.Item(15,16).Value = 14400
$ws.Cells.Item(15,17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(15,18).Value = 'Región del Maule'
$ws.Cells.Item(15,19).Value = 800
$ws.Cells.Item(15,20).Value = 18
$ws.Cells.Item(15,4).NumberFormat = $ws.Cells.Item(13,4).NumberFormat


Write-Host "Edit complete: dimension now includes rows 2-15"
